# Bump the "Förändrad" (Changed) date in column C for all data rows
# (rows 2-28) by one day: serial date 45524 -> 45525.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45524) {
        $cell.Value2 = 45525
    }
}
